# Report completed till index
#
# The "Ember.js" paragraph used to carry its own (redundant) run-level
# formatting, the blank paragraph right after it still carried a stray
# <w:color>, and "System architecture" carried a duplicated pPr/rPr block.
# This edit strips that redundant/duplicated run & paragraph formatting and
# appends two new (plain) paragraphs - "About testing environment" and
# "Scope & limitations" - after "System architecture".

$d = $word.ActiveDocument
$wN = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

function Get-ParaIndexByText($doc, $targetText) {
    $n = $doc.Paragraphs.Count
    for ($i = 1; $i -le $n; $i++) {
        $t = $doc.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $targetText) {
            return $i
        }
    }
    return -1
}

$emberIdx = Get-ParaIndexByText $d "Ember.js"
$sysIdx = Get-ParaIndexByText $d "System architecture"
$emptyIdx = $emberIdx + 1

if ($emberIdx -eq -1 -or $sysIdx -eq -1) {
    throw "Could not locate anchor paragraphs (Ember.js / System architecture)"
}

# Work bottom-to-top so the indices found above stay valid while we edit.

# 3) "System architecture" paragraph: drop its pPr/run rPr and add the two
#    new plain paragraphs right after it.
$pSys = $d.Paragraphs.Item($sysIdx)
$xmlSys = @"
<w:p xmlns:w="$wN"><w:r><w:t>System architecture</w:t></w:r></w:p><w:p xmlns:w="$wN"><w:r><w:t>About testing environment</w:t></w:r></w:p><w:p xmlns:w="$wN"><w:r><w:t>Scope &amp; limitations</w:t></w:r></w:p>
"@
$pSys.Range.InsertXML($xmlSys) | Out-Null

# 2) The blank paragraph right after "Ember.js": keep its pPr but drop the
#    stray <w:color>.
$pEmpty = $d.Paragraphs.Item($emptyIdx)
$xmlEmpty = @"
<w:p xmlns:w="$wN"><w:pPr><w:rPr><w:rFonts w:cs="Times New Roman"/><w:bCs/><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr></w:pPr></w:p>
"@
$pEmpty.Range.InsertXML($xmlEmpty) | Out-Null

# 1) "Ember.js" paragraph: keep its own pPr, strip the run's direct rPr.
$pEmber = $d.Paragraphs.Item($emberIdx)
$xmlEmber = @"
<w:p xmlns:w="$wN"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="222222"/><w:sz w:val="45"/><w:szCs w:val="45"/><w:lang w:val="en-IE" w:eastAsia="en-IN"/></w:rPr></w:pPr><w:r><w:t>Ember.js</w:t></w:r></w:p>
"@
$pEmber.Range.InsertXML($xmlEmber) | Out-Null
